# Applies the "Add more task" edit to Product Backlog.xlsx
#  - Project Backlog sheet (sheet1): row 4 & 5 Value changes, new Sprint values,
#    a new row 18 (No. 15), column C re-fit, table/autofilter grown to G18.
#  - Impediment Backlog sheet (sheet2): a new impediment row (No. 3) appended,
#    column C re-fit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Project Backlog"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Project Backlog")

# Update existing data: rows 4 & 5 get a higher Value and a Sprint number.
$ws1.Range("E4").Value = 200
$ws1.Range("G4").Value = 2

$ws1.Range("E5").Value = 200
$ws1.Range("G5").Value = 2

# Append a new backlog item as row 18 (No. 15).
$ws1.Range("A18").Value = 15

# The table/autofilter range grows by one row to include the new item.
$wb.Worksheets.Item("Project Backlog").ListObjects.Item("Table1").Resize($ws1.Range("A2:G18"))

# Column C is now best-fit (narrower, auto-sized) instead of a fixed width.
$ws1.Columns.Item(3).AutoFit()

# Restore the user's scroll position / selection as captured in the file.
$ws1.Application.ActiveWindow.ScrollRow = 2
$ws1.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Impediment Backlog"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Impediment Backlog")

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "Not solve"
$ws2.Range("C5").Value = "When finish task, not commit to SVN for other member verify it"

# Copy the red-font style used by the other "Status" cells in this column.
$ws2.Range("B3").Copy()
$ws2.Range("B5").PasteSpecial(-4122)  # xlPasteFormats

# Column C is now best-fit (wider, auto-sized) instead of a fixed width.
$ws2.Columns.Item(3).AutoFit()

$ws2.Range("C6").Select()

$wb.Save()
